$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 74.818184
$ws.Range("I6").Value = 42.3
$ws.Range("K6").Value = 126.9
$ws.Range("M6").Value = -14.89999999999999
$ws.Range("H8").Value = 162.5
$ws.Range("I8").Value = 162.5
$ws.Range("K8").Value = 487.5
$ws.Range("M8").Value = -348.5
$ws.Range("H12").Value = 396.375
$ws.Range("I12").Value = 396.375
$ws.Range("K12").Value = 396.375
$ws.Range("M12").Value = -226.375
$ws.Range("H18").Value = 998.38464
$ws.Range("I18").Value = 998.38464
$ws.Range("K18").Value = 998.38464
$ws.Range("M18").Value = -714.38464
$ws.Range("H40").Value = 4629.2856
$ws.Range("J40").Value = 6651.5
$ws.Range("L40").Value = 6651.5
$ws.Range("N40").Value = -7001.5
$ws.Range("H43").Value = 4704.364
$ws.Range("I43").Value = 4593.5
$ws.Range("K43").Value = 4593.5
$ws.Range("M43").Value = -4524.5
$ws.Range("H116").Value = 5126.5
$ws.Range("I116").Value = 3750
$ws.Range("K116").Value = 3750
$ws.Range("M116").Value = -308
$ws.Range("H135").Value = 706.75
$ws.Range("I135").Value = 723.6
$ws.Range("J135").Value = 678.6667
$ws.Range("K135").Value = 6512.400000000001
$ws.Range("L135").Value = 6108.0003
$ws.Range("M135").Value = -3977.400000000001
$ws.Range("N135").Value = -11178.0003
$ws.Range("H137").Value = 3049.3
$ws.Range("I137").Value = 1642.1428
$ws.Range("K137").Value = 4926.428400000001
$ws.Range("M137").Value = -2376.428400000001
$ws.Range("H138").Value = 2376.2
$ws.Range("I138").Value = 1971.9231
$ws.Range("K138").Value = 5915.7693
$ws.Range("M138").Value = -775.7692999999999

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 935.1667
$ws.Range("I2").Value = 935.1667
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 935.1667
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -822.1667
$ws.Range("N2").Value = $null
$ws.Range("H74").Value = 2847.7273
$ws.Range("I74").Value = 3281.9375
$ws.Range("K74").Value = 3281.9375
$ws.Range("M74").Value = -2407.9375
$ws.Range("H77").Value = 2847.7273
$ws.Range("I77").Value = 3281.9375
$ws.Range("K77").Value = 16409.6875
$ws.Range("M77").Value = -12041.6875
$ws.Range("H110").Value = 7620.3076
$ws.Range("I110").Value = 7620.3076
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 7620.3076
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -5575.3076
$ws.Range("N110").Value = $null
$ws.Range("H116").Value = 935.1667
$ws.Range("I116").Value = 935.1667
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 935.1667
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1358.8333
$ws.Range("N116").Value = $null
$ws.Range("H119").Value = 45000
$ws.Range("J119").Value = 45000
$ws.Range("L119").Value = 45000
$ws.Range("N119").Value = -54676
$ws.Range("H122").Value = 956.5
$ws.Range("I122").Value = 956.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2869.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -419.5
$ws.Range("N122").Value = $null
$ws.Range("H132").Value = 4442.8945
$ws.Range("I132").Value = 4442.8945
$ws.Range("K132").Value = 13328.6835
$ws.Range("M132").Value = -10798.6835

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 935.1667
$ws.Range("I3").Value = 935.1667
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 935.1667
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -821.1667
$ws.Range("N3").Value = $null
$ws.Range("H107").Value = 939.8333
$ws.Range("I107").Value = 887.8
$ws.Range("K107").Value = 887.8
$ws.Range("M107").Value = 1032.2

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1211
$ws.Range("I31").Value = 1141.4166
$ws.Range("K31").Value = 1141.4166
$ws.Range("M31").Value = -846.4166
$ws.Range("H34").Value = 1211
$ws.Range("I34").Value = 1141.4166
$ws.Range("K34").Value = 1141.4166
$ws.Range("M34").Value = -939.4166
$ws.Range("H58").Value = 1827.5
$ws.Range("I58").Value = 1853.3334
$ws.Range("K58").Value = 1853.3334
$ws.Range("M58").Value = -1650.3334
$ws.Range("H99").Value = 1456.4445
$ws.Range("I99").Value = 1259.8334
$ws.Range("K99").Value = 1259.8334
$ws.Range("M99").Value = 238.1666
$ws.Range("H126").Value = 1456.4445
$ws.Range("I126").Value = 1259.8334
$ws.Range("K126").Value = 3779.5002
$ws.Range("M126").Value = -1309.5002
$ws.Range("H132").Value = 2809.25
$ws.Range("I132").Value = 2712.3333
$ws.Range("K132").Value = 8136.999899999999
$ws.Range("M132").Value = -5606.999899999999
$ws.Range("H135").Value = 135390
$ws.Range("J135").Value = 135390
$ws.Range("L135").Value = 135390
$ws.Range("N135").Value = -145530
$ws.Range("H136").Value = 1827.5
$ws.Range("I136").Value = 1853.3334
$ws.Range("K136").Value = 5560.0002
$ws.Range("M136").Value = -3010.0002

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1500
$ws.Range("I64").Value = 1500
$ws.Range("K64").Value = 4500
$ws.Range("M64").Value = -4230
$ws.Range("H67").Value = 1500
$ws.Range("I67").Value = 1500
$ws.Range("K67").Value = 4500
$ws.Range("M67").Value = -3564
$ws.Range("H68").Value = 2681.5862
$ws.Range("I68").Value = 2113.7144
$ws.Range("K68").Value = 6341.1432
$ws.Range("M68").Value = -5530.1432
$ws.Range("H71").Value = 2681.5862
$ws.Range("I71").Value = 2113.7144
$ws.Range("K71").Value = 19023.4296
$ws.Range("M71").Value = -14967.4296
$ws.Range("H106").Value = 7307.5
$ws.Range("J106").Value = 8500
$ws.Range("L106").Value = 25500
$ws.Range("N106").Value = -27392
$ws.Range("H114").Value = 1915.3334
$ws.Range("I114").Value = 1798.4
$ws.Range("J114").Value = 2500
$ws.Range("K114").Value = 5395.200000000001
$ws.Range("L114").Value = 7500
$ws.Range("M114").Value = -2141.200000000001
$ws.Range("N114").Value = -14008

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = $null
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = $null
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = $null
$ws.Range("H102").Value = 2186.8572
$ws.Range("I102").Value = 2186.8572
$ws.Range("K102").Value = 2186.8572
$ws.Range("M102").Value = -564.8571999999999
$ws.Range("H132").Value = 2377.2942
$ws.Range("I132").Value = 2481
$ws.Range("K132").Value = 7443
$ws.Range("M132").Value = -4913
$ws.Range("H139").Value = 21995
$ws.Range("J139").Value = 21995
$ws.Range("L139").Value = 21995
$ws.Range("N139").Value = -32275

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2498.8333
$ws.Range("I68").Value = 2498.8333
$ws.Range("K68").Value = 2498.8333
$ws.Range("M68").Value = -1749.8333
$ws.Range("H71").Value = 2498.8333
$ws.Range("I71").Value = 2498.8333
$ws.Range("K71").Value = 12494.1665
$ws.Range("M71").Value = -8750.166499999999
$ws.Range("H93").Value = 777.3333
$ws.Range("I93").Value = 749.5
$ws.Range("K93").Value = 749.5
$ws.Range("M93").Value = 498.5
$ws.Range("H122").Value = 33800.2
$ws.Range("I122").Value = 37502
$ws.Range("J122").Value = 31332.334
$ws.Range("K122").Value = 112506
$ws.Range("L122").Value = 93997.00199999999
$ws.Range("M122").Value = -110056
$ws.Range("N122").Value = -98897.00199999999
$ws.Range("H128").Value = 60214.5
$ws.Range("J128").Value = 60214.5
$ws.Range("L128").Value = 60214.5
$ws.Range("N128").Value = -70174.5
$ws.Range("H130").Value = 79672.5
$ws.Range("J130").Value = 79672.5
$ws.Range("L130").Value = 79672.5
$ws.Range("N130").Value = -89712.5
$ws.Range("H132").Value = 3413.7144
$ws.Range("I132").Value = 2649.5
$ws.Range("K132").Value = 7948.5
$ws.Range("M132").Value = -5418.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2535.75
$ws.Range("I122").Value = 2535.75
$ws.Range("K122").Value = 7607.25
$ws.Range("M122").Value = -5157.25
$ws.Range("H126").Value = 1991.8572
$ws.Range("I126").Value = 1992
$ws.Range("K126").Value = 5976
$ws.Range("M126").Value = -3506
$ws.Range("H132").Value = 7273
$ws.Range("I132").Value = 5581.2
$ws.Range("K132").Value = 16743.6
$ws.Range("M132").Value = -14213.6
